{"js": "// The document has a paragraph holding a single Word field whose field\n// code is the M2Doc expression ` m:'Mona_Lisa.jpg'.asImage().rotate(90) `\n// (classic `{ FIELD }` begin/instrText/end run structure, with the\n// `'Mona_Lisa.jpg'.asImage().rotate(90)` part colored orange).\n//\n// The commit replaces that field with plain literal text runs spelling\n// out the same characters as an M2Doc template tag: the leading\n// `fldChar begin` + leading space collapse into a literal \"{\" run, every\n// `w:instrText` run becomes an equivalent `w:t` run (keeping its\n// formatting/rPr untouched), and the trailing space + `fldChar end`\n// collapse into a literal \"}\" run. The bookmark in the middle of the\n// field is left exactly where it was.\n\n// Find the paragraph that owns the field (robust to index drift).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (const p of paragraphs.items) {\n  const fields = p.fields;\n  fields.load(\"items\");\n  await context.sync();\n  if (fields.items.length > 0) {\n    targetParagraph = p;\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not find the paragraph containing the field.\");\n}\n\n// Rebuild that paragraph's OOXML: same runs/rPr as the field's instrText\n// runs, but using <w:t> instead of <w:fldChar>/<w:instrText>, with the\n// field delimiters turned into literal \"{\" and \"}\" runs.\nconst color = '<w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr>';\nconst newParagraphXml =\n  '<w:p w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" w:rsidP=\"00F5495F\">' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r w:rsidR=\"00DE6D5A\"><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    '<w:r w:rsidR=\"004B598D\">' + color + \"<w:t>'</w:t></w:r>\" +\n    '<w:r w:rsidR=\"003C367E\" w:rsidRPr=\"003C367E\">' + color + '<w:t>Mona_Lisa</w:t></w:r>' +\n    '<w:r w:rsidR=\"00882765\">' + color + '<w:t>.jpg</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r w:rsidR=\"004B598D\">' + color + \"<w:t>'.asImage()</w:t></w:r>\" +\n    '<w:r w:rsidR=\"00484D7C\">' + color + '<w:t>.</w:t></w:r>' +\n    '<w:r w:rsidR=\"00215769\">' + color + '<w:t>rotate</w:t></w:r>' +\n    '<w:r w:rsidR=\"00484D7C\">' + color + '<w:t>(</w:t></w:r>' +\n    '<w:r w:rsidR=\"00215769\">' + color + '<w:t>90</w:t></w:r>' +\n    '<w:r w:rsidR=\"00FB48D7\">' + color + '<w:t>)</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>';\n\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + newParagraphXml + '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\ntargetParagraph.insertOoxml(ooxmlPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document has a paragraph holding a single Word field whose field\n# code is the M2Doc expression ` m:'Mona_Lisa.jpg'.asImage().rotate(90) `\n# (classic `{ FIELD }` begin/instrText/end run structure, with the\n# `'Mona_Lisa.jpg'.asImage().rotate(90)` part colored orange).\n#\n# This rewrites that field into plain literal text runs spelling out the\n# same characters as an M2Doc template tag: the leading `fldChar begin` +\n# leading space collapse into a literal \"{\" run, every `w:instrText` run\n# becomes an equivalent `w:t` run (keeping its formatting/rPr untouched),\n# and the trailing space + `fldChar end` collapse into a literal \"}\" run.\n# The bookmark in the middle of the field is left exactly where it was.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph that owns the field (robust to index drift).\n$targetParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Fields.Count -gt 0) {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($targetParagraph -eq $null) {\n    throw \"Could not find the paragraph containing the field.\"\n}\n\n$newParagraphXml = '<w:p w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" w:rsidP=\"00F5495F\">' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r w:rsidR=\"00DE6D5A\"><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    '<w:r w:rsidR=\"004B598D\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>' + [char]39 + '</w:t></w:r>' +\n    '<w:r w:rsidR=\"003C367E\" w:rsidRPr=\"003C367E\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>Mona_Lisa</w:t></w:r>' +\n    '<w:r w:rsidR=\"00882765\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>.jpg</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r w:rsidR=\"004B598D\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>' + [char]39 + '.asImage()</w:t></w:r>' +\n    '<w:r w:rsidR=\"00484D7C\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>.</w:t></w:r>' +\n    '<w:r w:rsidR=\"00215769\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>rotate</w:t></w:r>' +\n    '<w:r w:rsidR=\"00484D7C\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>(</w:t></w:r>' +\n    '<w:r w:rsidR=\"00215769\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>90</w:t></w:r>' +\n    '<w:r w:rsidR=\"00FB48D7\"><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>)</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n    '</w:p>'\n\n$ooxmlPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n            '<pkg:xmlData>' +\n                '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n                    '<w:body>' + $newParagraphXml + '</w:body>' +\n                '</w:document>' +\n            '</pkg:xmlData>' +\n        '</pkg:part>' +\n    '</pkg:package>'\n\n$targetParagraph.Range.InsertXML($ooxmlPackage)\n"}
